$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Activate()

# Row 44
$ws.Range("A44").Value = 35
$ws.Range("B44").Value = "ActualFilingDate"
$ws.Range("C44").Value = "實際報送日期"
$ws.Range("D44").Value = "Decimald"
$ws.Range("E44").Value = 8

# Row 45
$ws.Range("A45").Value = 36
$ws.Range("B45").Value = "ActualFilingMark"
$ws.Range("C45").Value = "實際報送記號"
$ws.Range("D45").Value = "VARCHAR2"
$ws.Range("E45").Value = 3

# Update selection to match the new active range
$ws.Range("B44:E45").Select()
